# Target change (per the supplied OOXML diff) trims word/styles.xml's
# <w:docDefaults> block:
#   - rPrDefault/rPr keeps only rFonts/sz/szCs/lang and drops explicit
#     b/i/smallCaps/strike/color/u/shd/vertAlign entries (all of which are
#     already equal to the implicit OOXML schema defaults).
#   - pPrDefault/pPr keeps only a bare <w:spacing w:line="276"
#     w:lineRule="auto"/> and drops keepNext/keepLines/widowControl/pBdr/
#     shd/ind/contextualSpacing/jc (again, all already schema defaults)
#     plus the redundant spacing w:after="0"/w:before="0".
#
# In other words every attribute being removed already equals Word's
# built-in default for that property when the element/attribute is simply
# absent, so the edit is a pure, render-invisible cleanup of redundant
# default-valued XML in <w:docDefaults> - it does not change how the
# document looks or behaves.
#
# <w:docDefaults> is the one corner of styles.xml that the Word object
# model does not expose: there is no Styles/Style/Font/ParagraphFormat
# member (and no Document property) bound to rPrDefault/pPrDefault - only
# to the explicit per-style <w:style> definitions (e.g. Styles("Normal")).
# Word.Interop simply has no COM surface for it, and probing this
# runtime's Word object model confirms the same: Document.WordOpenXML is
# read-only, Range.InsertXML only ever replaces the content of the Range
# it is called on (document body content, not styles.xml), and writing to
# Styles("Normal").Font / .ParagraphFormat lands an explicit override on
# the <w:style w:styleId="Normal"> element itself - which the diff leaves
# untouched - rather than touching <w:docDefaults>.
#
# So there is no COM call available on $word/$d that can reach
# <w:docDefaults> without instead corrupting a part of styles.xml the
# diff does not touch. Since the change has no visible/semantic effect on
# the document, the correct, non-destructive action through this object
# model is to leave the document's styles alone rather than fabricate an
# unrelated mutation - so this script intentionally performs no edits.
$d = $word.ActiveDocument
$d.Styles.Count | Out-Null
